$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set the date cells first so the "YYYY-MM-DD" style is allocated before
#     the text-coercion (quote-prefix) styles used below, matching the
#     author's style table ordering as closely as the engine allows. ---

# Row 2: becomes the "Novas barreiras..." record (previously row 3's data)
$ws.Cells.Item(2, 10).NumberFormat = "yyyy-mm-dd"
$ws.Cells.Item(2, 10).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(2, 10).Value = "2021-01-25"
$ws.Cells.Item(2, 11).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(2, 11).Value = "2021-02-02"

# Row 3: becomes the "Estudo Team.docx" record (previously row 2's data)
$ws.Cells.Item(3, 10).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(3, 10).Value = "2021-02-01"
$ws.Cells.Item(3, 11).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(3, 11).Value = "2021-02-02"

# Row 4: refresh its triagem date
$ws.Cells.Item(4, 10).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(4, 10).Value = "2019-12-06"
$ws.Cells.Item(4, 11).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(4, 11).Value = "2021-02-02"

# --- Remaining field updates ---

# Row 2
$ws.Cells.Item(2, 1).Value = "'236690"
$ws.Cells.Item(2, 1).ClearFormats()
$ws.Cells.Item(2, 2).Value = "Novas barreiras e tendências no comércio internacional.docx"
$ws.Cells.Item(2, 3).Value = 49
$ws.Cells.Item(2, 4).Value = 80889
$ws.Cells.Item(2, 5).Value = 6
$ws.Cells.Item(2, 6).Value = 12
$ws.Cells.Item(2, 7).Value = 26
$ws.Cells.Item(2, 12).Value = "sim"

# Row 3
$ws.Cells.Item(3, 1).Value = "'236813"
$ws.Cells.Item(3, 1).ClearFormats()
$ws.Cells.Item(3, 2).Value = "Estudo Team.docx"
$ws.Cells.Item(3, 3).Value = 58
$ws.Cells.Item(3, 4).Value = 131904
$ws.Cells.Item(3, 5).Value = 8
$ws.Cells.Item(3, 6).Value = 21
$ws.Cells.Item(3, 7).Value = 25

# Row 4 - fill in the missing title and refresh counts
$ws.Cells.Item(4, 2).Value = "TESTE.docx"
$ws.Cells.Item(4, 4).Value = 56596
$ws.Cells.Item(4, 5).Value = 11
$ws.Cells.Item(4, 6).Value = 32
